$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first worksheet "3rd" -> "Order"
# ---------------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item(1)
$wsOrder.Name = "Order"

# ---------------------------------------------------------------------------
# 2. Fix up the defined names (Print_Area / Print_Titles) that referenced the
#    old sheet name '3rd' so they point at the renamed sheet "Order".
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*'3rd'*") {
        $n.RefersTo = ($n.RefersTo -replace "'3rd'", "Order")
    }
}

# ---------------------------------------------------------------------------
# 3. Re-word the purchase-purpose sentence in A6 to use generic placeholders
#    instead of the hard coded "20년 9월 / 2020" date.
# ---------------------------------------------------------------------------
$wsOrder.Range("A6").Value = "1. 구입 목적 (Mục đích mua) : ….년 ... 월 생산계획에 따른 Consumable part 구입 품의 (Đề nghị mua vật tư tiêu hao dựa trên kế hoạch sản xuất tháng …. năm ….)"

# ---------------------------------------------------------------------------
# 4. Re-enter the "Tên tiếng anh/ tiếng hàn" header (E11) - same visible text
#    but split across four runs (바탕 / Times New Roman alternating).
# ---------------------------------------------------------------------------
$cellE11 = $wsOrder.Range("E11")
$cellE11.Value = "
Tên tiếng anh/ tiếng hàn"

$run1 = $cellE11.Characters(1, 1)
$run1.Font.Name = "바탕"
$run1.Font.Bold = $true
$run1.Font.Size = 24

$run2 = $cellE11.Characters(2, 3)
$run2.Font.Name = "Times New Roman"
$run2.Font.Bold = $true
$run2.Font.Size = 24

$run3 = $cellE11.Characters(5, 1)
$run3.Font.Name = "바탕"
$run3.Font.Bold = $true
$run3.Font.Size = 24

$run4 = $cellE11.Characters(6, 20)
$run4.Font.Name = "Times New Roman"
$run4.Font.Bold = $true
$run4.Font.Size = 24

# ---------------------------------------------------------------------------
# 5. Move the active selection to F8 (matches the saved cursor position).
# ---------------------------------------------------------------------------
$wsOrder.Activate()
$wsOrder.Range("F8").Select()
